$d = $word.ActiveDocument

# The document has two placeholder inline pictures (both pointing at the
# same 1x1 embedded png) that are being replaced by plain-text hyperlinks
# pointing at the real, hosted images on ura.gov.sg.

$url1 = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Industrial/Special-Control-Area-1.jpg?h=416&w=750"
$url2 = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Industrial/Special-Control-Area-2.jpg?h=383&w=750"

# First picture -> first hyperlink.
$shape1 = $d.InlineShapes.Item(1)
$range1 = $shape1.Range
$shape1.Delete()
$d.Hyperlinks.Add($range1, $url1)

# Second picture -> second hyperlink.
$shape2 = $d.InlineShapes.Item(1)
$range2 = $shape2.Range
$shape2.Delete()
$d.Hyperlinks.Add($range2, $url2)
